# Update need_to_buy.xlsx data: refresh the daily table (A2:F15) with the
# latest values - the series rolled forward by one day (oldest day dropped,
# a new day appended) and the forecast/buy columns were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45966, 8615.14871708983, 8324.21427667692, 12964, 4195.66946498772, -18.504844097307),
    @(45967, 8642.25204484622, 8394.84851461668, 9004, 4651.90447591225, 168.448041272039),
    @(45968, 8522.36699430934, 7595.0711559492, 9004, 4609.4439440977, 133.354795835287),
    @(45969, 2959.12979841581, 5162.82776751609, 9004, 4148.59552844887, 12.80930399854),
    @(45970, 2676.75776179279, 5038.18361509195, 9004, 3968.35120640539, 0.105617562388791),
    @(45971, 9273.47168962819, 8557.86312822128, 9004, 5243.44741970636, 199.887939496985),
    @(45972, 9273.47299076555, 8929.25448444316, 9004, 5243.45118205216, 215.362736103972),
    @(45973, 9273.47299076555, 8634.72379961751, 9004, 5243.45118205216, 203.090624236236),
    @(45974, 9273.47299076555, 7980.87006984986, 9004, 5243.45118205216, 175.846718829251),
    @(45975, 9273.47299076555, 7042.05853267733, 9004, 5243.45118205216, 136.729571447062),
    @(45976, 3551.20396235217, 5004.12477501231, 9004, 4928.56350043268, 38.6953448102079),
    @(45977, 3412.56435984501, 4869.63736200861, 9004, 4921.09168721732, 32.7803770510806),
    @(45978, 10116.508538002, 8420.4437654909, 9004, 5829.29625007391, 218.572500648534),
    @(45979, 10116.508538002, 8768.45052111632, 9004, 5829.29625007391, 233.072782132926)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
